$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.836.71"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.637.94"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'216.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "'0.5076"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'0.2585"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").Value = "'0.06445"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "'19.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").Value = "'0.07796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "'4.288"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.864.83"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "1.635.73"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "'0.5633"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").Value = "0.0₅7609"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "'63.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "25.873.74"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "'195.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'4.333"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("D22").Value = "'9.899"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").Value = "'6.110"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'1.785"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.99%  "
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").Value = "'139.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").Value = "'6.796"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("D29").Value = "'15.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "'1.243"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "'0.04886"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "'3.303"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "'3.235"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").Value = "'1.558"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").Value = "'2.369"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "'0.9045"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D38").Value = "1.129.89"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").Value = "'0.5509"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "'0.9946"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").Value = "'5.535"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "'0.8018"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "'97.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "1.775.20"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").Value = "  -6.88%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.4445"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'55.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").Value = "'7.714"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.32%  "
$ws.Range("D50").Value = "'0.05052"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").Value = "'1.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
